$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string into a cell, forcing text storage when the
# string would otherwise be auto-coerced into a number (e.g. "210.89",
# "1.00", "0.0610") by Excel's normal type inference on a Range.Value write.
# Temporarily switching the cell to a text NumberFormat for the assignment,
# then clearing formats again, keeps the value exactly as typed without
# leaving any formatting change behind on the cell.
function Set-TextCell($range, [string]$text) {
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.ClearFormats()
    } else {
        $range.Value = $text
    }
}


Set-TextCell $ws.Range('D2') '27.461.29'
Set-TextCell $ws.Range('E2') '  -1.21%  '

Set-TextCell $ws.Range('D3') '1.614.07'
Set-TextCell $ws.Range('E3') '  -2.09%  '

Set-TextCell $ws.Range('D5') '210.89'
Set-TextCell $ws.Range('E5') '  -1.18%  '

Set-TextCell $ws.Range('E6') '  -1.58%  '

Set-TextCell $ws.Range('E7') '  +0.04%  '

Set-TextCell $ws.Range('D8') '22.73'
Set-TextCell $ws.Range('E8') '  -1.79%  '

Set-TextCell $ws.Range('E9') '  +0.50%  '

Set-TextCell $ws.Range('D10') '0.0610'
Set-TextCell $ws.Range('E10') '  -0.80%  '

Set-TextCell $ws.Range('D11') '0.0886'
Set-TextCell $ws.Range('E11') '  -0.40%  '

Set-TextCell $ws.Range('D12') '1.842.73'
Set-TextCell $ws.Range('E12') '  -2.14%  '

Set-TextCell $ws.Range('D13') '1.612.55'
Set-TextCell $ws.Range('E13') '  -2.17%  '

Set-TextCell $ws.Range('E14') '  -0.53%  '

Set-TextCell $ws.Range('D15') '0.546'
Set-TextCell $ws.Range('E15') '  -3.25%  '

Set-TextCell $ws.Range('D16') '64.93'
Set-TextCell $ws.Range('E16') '  +0.73%  '

Set-TextCell $ws.Range('D17') '27.426.12'
Set-TextCell $ws.Range('E17') '  -1.27%  '

Set-TextCell $ws.Range('D18') '231.81'
Set-TextCell $ws.Range('E18') '  -0.80%  '

Set-TextCell $ws.Range('E19') '  -1.49%  '

Set-TextCell $ws.Range('D20') '7.49'
Set-TextCell $ws.Range('E20') '  -2.48%  '

Set-TextCell $ws.Range('E21') '  +0.13%  '

Set-TextCell $ws.Range('E22') '  +0.16%  '

Set-TextCell $ws.Range('D23') '10.15'
Set-TextCell $ws.Range('E23') '  +0.46%  '

Set-TextCell $ws.Range('E24') '  +5.32%  '

Set-TextCell $ws.Range('D25') '149.84'
Set-TextCell $ws.Range('E25') '  -0.45%  '

Set-TextCell $ws.Range('E26') '  -1.78%  '

Set-TextCell $ws.Range('B27') 'BinanceUSD'
Set-TextCell $ws.Range('C27') 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell $ws.Range('D27') '1.00'
Set-TextCell $ws.Range('E27') '  +0.10%  '

Set-TextCell $ws.Range('B28') 'Stellar'
Set-TextCell $ws.Range('C28') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws.Range('D28') '0.111'
Set-TextCell $ws.Range('E28') '  -1.45%  '

Set-TextCell $ws.Range('D29') '15.49'
Set-TextCell $ws.Range('E29') '  -1.08%  '

Set-TextCell $ws.Range('E30') '  -1.19%  '

Set-TextCell $ws.Range('E31') '  -1.40%  '

Set-TextCell $ws.Range('D32') '3.26'
Set-TextCell $ws.Range('E32') '  -1.63%  '

Set-TextCell $ws.Range('D33') '1.471.22'
Set-TextCell $ws.Range('E33') '  +1.92%  '

Set-TextCell $ws.Range('E34') '  -3.48%  '

Set-TextCell $ws.Range('D35') '1.53'
Set-TextCell $ws.Range('E35') '  -3.37%  '

Set-TextCell $ws.Range('D36') '0.972'
Set-TextCell $ws.Range('E36') '  +11.69%  '

Set-TextCell $ws.Range('E37') '  -0.64%  '

Set-TextCell $ws.Range('E38') '  -0.66%  '

Set-TextCell $ws.Range('D39') '0.554'
Set-TextCell $ws.Range('E39') '  -2.80%  '

Set-TextCell $ws.Range('D40') '0.859'
Set-TextCell $ws.Range('E40') '  -2.83%  '

Set-TextCell $ws.Range('D42') '66.91'
Set-TextCell $ws.Range('E42') '  +0.33%  '

Set-TextCell $ws.Range('D43') '0.986'
Set-TextCell $ws.Range('E43') '  -4.61%  '

Set-TextCell $ws.Range('D44') '2.45'
Set-TextCell $ws.Range('E44') '  -0.47%  '

Set-TextCell $ws.Range('E45') '  -2.74%  '

Set-TextCell $ws.Range('D46') '1.753.90'
Set-TextCell $ws.Range('E46') '  -2.15%  '

Set-TextCell $ws.Range('D47') '5.21'
Set-TextCell $ws.Range('E47') '  -6.88%  '

Set-TextCell $ws.Range('E48') '  +0.28%  '

Set-TextCell $ws.Range('D49') '86.55'
Set-TextCell $ws.Range('E49') '  +0.23%  '

Set-TextCell $ws.Range('E50') '  -2.00%  '

Set-TextCell $ws.Range('D51') '0.100'
Set-TextCell $ws.Range('E51') '  +0.89%  '
